# Re-executes the RAD UI test suite: updates the Date (and, where the
# run result changed, Result) column for every already-recorded test row
# across the Estimated / Existing / NewTaxReturn / Personal_EL /
# Personal_IND / Personal_JNT sheets, in the order the runs completed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Tue Jan 30 22:39:43 EST 2024"
$ws.Range("B3").Value = "Tue Jan 30 22:40:38 EST 2024"
$ws.Range("B4").Value = "Tue Jan 30 22:41:28 EST 2024"
$ws.Range("B5").Value = "Tue Jan 30 22:42:18 EST 2024"
$ws.Range("B6").Value = "Tue Jan 30 22:43:08 EST 2024"
$ws.Range("B7").Value = "Tue Jan 30 22:43:57 EST 2024"
$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Tue Jan 30 22:44:47 EST 2024"
$ws.Range("B3").Value = "Tue Jan 30 22:45:37 EST 2024"
$ws.Range("B4").Value = "Tue Jan 30 22:46:26 EST 2024"
$ws.Range("B5").Value = "Tue Jan 30 22:47:16 EST 2024"
$ws.Range("B6").Value = "Tue Jan 30 22:48:05 EST 2024"
$ws.Range("B7").Value = "Tue Jan 30 22:48:55 EST 2024"
$ws.Range("B8").Value = "Tue Jan 30 22:49:44 EST 2024"
$ws.Range("B9").Value = "Tue Jan 30 22:50:35 EST 2024"
$ws.Range("B10").Value = "Tue Jan 30 22:51:26 EST 2024"
$ws.Range("B11").Value = "Tue Jan 30 22:52:16 EST 2024"
$ws.Range("B12").Value = "Tue Jan 30 22:53:06 EST 2024"
$ws.Range("B13").Value = "Tue Jan 30 22:53:55 EST 2024"
$ws.Range("B14").Value = "Tue Jan 30 22:54:44 EST 2024"
$ws.Range("A15").Value = "Fail"
$ws.Range("B15").Value = "Tue Jan 30 22:55:34 EST 2024"
$ws.Range("B16").Value = "Tue Jan 30 22:57:06 EST 2024"
$ws.Range("B17").Value = "Tue Jan 30 22:57:56 EST 2024"
$ws.Range("B18").Value = "Tue Jan 30 22:58:46 EST 2024"
$ws.Range("B19").Value = "Tue Jan 30 22:59:35 EST 2024"
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Tue Jan 30 23:00:25 EST 2024"
$ws.Range("B3").Value = "Tue Jan 30 23:01:14 EST 2024"
$ws.Range("B4").Value = "Tue Jan 30 23:02:03 EST 2024"
$ws.Range("B5").Value = "Tue Jan 30 23:02:52 EST 2024"
$ws.Range("B6").Value = "Tue Jan 30 23:03:40 EST 2024"
$ws.Range("B7").Value = "Tue Jan 30 23:04:29 EST 2024"
$ws.Range("B8").Value = "Tue Jan 30 23:05:18 EST 2024"
$ws.Range("B9").Value = "Tue Jan 30 23:06:07 EST 2024"
$ws.Range("B10").Value = "Tue Jan 30 23:06:56 EST 2024"
$ws.Range("B11").Value = "Tue Jan 30 23:07:46 EST 2024"
$ws.Range("B12").Value = "Tue Jan 30 23:08:35 EST 2024"
$ws.Range("B13").Value = "Tue Jan 30 23:09:24 EST 2024"
$ws.Range("B14").Value = "Tue Jan 30 23:10:13 EST 2024"
$ws.Range("B15").Value = "Tue Jan 30 23:11:02 EST 2024"
$ws.Range("B16").Value = "Tue Jan 30 23:11:51 EST 2024"
$ws.Range("B17").Value = "Tue Jan 30 23:12:40 EST 2024"
$ws.Range("B18").Value = "Tue Jan 30 23:13:29 EST 2024"
$ws.Range("B19").Value = "Tue Jan 30 23:14:18 EST 2024"
$ws.Range("B20").Value = "Tue Jan 30 23:15:08 EST 2024"
$ws.Range("B21").Value = "Tue Jan 30 23:15:56 EST 2024"
$ws.Range("B22").Value = "Tue Jan 30 23:16:45 EST 2024"
$ws.Range("B23").Value = "Tue Jan 30 23:17:34 EST 2024"
$ws.Range("B24").Value = "Tue Jan 30 23:18:23 EST 2024"
$ws.Range("B25").Value = "Tue Jan 30 23:19:11 EST 2024"
$ws.Range("B26").Value = "Tue Jan 30 23:20:00 EST 2024"
$ws.Range("B27").Value = "Tue Jan 30 23:20:50 EST 2024"
$ws.Range("B28").Value = "Tue Jan 30 23:21:39 EST 2024"
$ws.Range("B29").Value = "Tue Jan 30 23:22:27 EST 2024"
$ws.Range("B30").Value = "Tue Jan 30 23:23:16 EST 2024"
$ws.Range("B31").Value = "Tue Jan 30 23:24:05 EST 2024"
$ws.Range("B32").Value = "Tue Jan 30 23:24:53 EST 2024"
$ws.Range("B33").Value = "Tue Jan 30 23:25:42 EST 2024"
$ws.Range("B34").Value = "Tue Jan 30 23:26:31 EST 2024"
$ws.Range("B35").Value = "Tue Jan 30 23:27:20 EST 2024"
$ws.Range("B36").Value = "Tue Jan 30 23:28:08 EST 2024"
$ws.Range("B37").Value = "Tue Jan 30 23:28:56 EST 2024"
$ws.Range("B38").Value = "Tue Jan 30 23:29:45 EST 2024"
$ws.Range("B39").Value = "Tue Jan 30 23:30:34 EST 2024"
$ws.Range("B40").Value = "Tue Jan 30 23:31:23 EST 2024"
$ws.Range("B41").Value = "Tue Jan 30 23:32:12 EST 2024"
$ws.Range("B42").Value = "Tue Jan 30 23:33:00 EST 2024"
$ws.Range("B43").Value = "Tue Jan 30 23:33:49 EST 2024"
$ws.Range("B44").Value = "Tue Jan 30 23:34:38 EST 2024"
$ws.Range("B45").Value = "Tue Jan 30 23:35:27 EST 2024"
$ws.Range("B46").Value = "Tue Jan 30 23:36:16 EST 2024"
$ws.Range("B47").Value = "Tue Jan 30 23:37:04 EST 2024"
$ws.Range("B48").Value = "Tue Jan 30 23:37:53 EST 2024"
$ws.Range("B49").Value = "Tue Jan 30 23:38:42 EST 2024"
$ws.Range("B50").Value = "Tue Jan 30 23:39:30 EST 2024"
$ws.Range("B51").Value = "Tue Jan 30 23:40:19 EST 2024"
$ws.Range("B52").Value = "Tue Jan 30 23:41:08 EST 2024"
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Jan 30 23:41:57 EST 2024"
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Jan 30 23:42:45 EST 2024"
$ws.Range("B4").Value = "Tue Jan 30 23:43:32 EST 2024"
$ws.Range("B5").Value = "Tue Jan 30 23:44:18 EST 2024"
$ws.Range("B6").Value = "Tue Jan 30 23:45:05 EST 2024"
$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Tue Jan 30 23:45:53 EST 2024"
$ws.Range("B4").Value = "Tue Jan 30 23:46:45 EST 2024"
$ws.Range("B5").Value = "Tue Jan 30 23:47:37 EST 2024"
$ws.Range("B6").Value = "Tue Jan 30 23:48:29 EST 2024"
